$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = [double]"21.13602999991246"
$ws.Range("C2").Value = [double]"381"
$ws.Range("E2").Value = [double]"-7.720552620793872E-07"
$ws.Range("F2").Value = [double]"0.2389096387727623"
$ws.Range("G2").Value = [double]"3669.052557726704"
$ws.Range("H2").Value = [double]"0.5760623394560491"
$ws.Range("B3").Value = [double]"21.91985209994788"
$ws.Range("C3").Value = [double]"5"
$ws.Range("F3").Value = [double]"0.2818152875142326"
$ws.Range("G3").Value = [double]"3800.132541317691"
$ws.Range("H3").Value = [double]"0.5768180941485583"
$ws.Range("B4").Value = [double]"22.70966652997307"
$ws.Range("E4").Value = [double]"3.860262620793873E-07"
$ws.Range("F4").Value = [double]"0.3014897732901751"
$ws.Range("G4").Value = [double]"3989.724155585123"
$ws.Range("H4").Value = [double]"0.5692039259952929"
$ws.Range("B5").Value = [double]"23.48373150006158"
$ws.Range("F5").Value = [double]"0.3136068395135485"
$ws.Range("G5").Value = [double]"4196.359215999511"
$ws.Range("H5").Value = [double]"0.5596215741141717"
$ws.Range("B6").Value = [double]"24.36210549998146"
$ws.Range("F6").Value = [double]"0.3226878663647829"
$ws.Range("G6").Value = [double]"4399.733934935221"
$ws.Range("H6").Value = [double]"0.5537176988485363"
$ws.Range("B7").Value = [double]"25.37542010995438"
$ws.Range("E7").Value = [double]"123.354268245"
$ws.Range("F7").Value = [double]"0.3202428064628157"
$ws.Range("G7").Value = [double]"4709.601962012134"
$ws.Range("H7").Value = [double]"0.5388017992737748"
$ws.Range("B8").Value = [double]"26.54542002997925"
$ws.Range("E8").Value = [double]"131.742805953"
$ws.Range("F8").Value = [double]"0.3646651455146169"
$ws.Range("G8").Value = [double]"4901.853013030957"
$ws.Range("H8").Value = [double]"0.5415384745199744"
$ws.Range("D9").Value = [double]"17"
$ws.Range("E9").Value = [double]"131.8719256233472"
$ws.Range("F9").Value = [double]"0.3734427761355451"
$ws.Range("G9").Value = [double]"5034.280137212208"
$ws.Range("H9").Value = [double]"0.5566181159215553"
$ws.Range("B10").Value = [double]"28.84944070996335"
$ws.Range("E10").Value = [double]"277.734877506"
$ws.Range("F10").Value = [double]"0.3675990765247539"
$ws.Range("G10").Value = [double]"5265.494010134662"
$ws.Range("H10").Value = [double]"0.547896183234392"
$ws.Range("B11").Value = [double]"29.67297186996335"
$ws.Range("C11").Value = [double]"0"
$ws.Range("D11").Value = [double]"16"
$ws.Range("E11").Value = [double]"155.9879130427576"
$ws.Range("F11").Value = [double]"0.3832457607533119"
$ws.Range("G11").Value = [double]"5400.305099148954"
$ws.Range("H11").Value = [double]"0.5494684341934603"
$ws.Range("B12").Value = [double]"30.56014719996283"
$ws.Range("C12").Value = [double]"19"
$ws.Range("D12").Value = [double]"31"
$ws.Range("E12").Value = [double]"355.4414327389016"
$ws.Range("F12").Value = [double]"0.3766049914101784"
$ws.Range("G12").Value = [double]"5704.430718030029"
$ws.Range("H12").Value = [double]"0.5357265029684941"
$ws.Range("B13").Value = [double]"31.47516354999119"
$ws.Range("C13").Value = [double]"14.00000001186676"
$ws.Range("D13").Value = [double]"50"
$ws.Range("E13").Value = [double]"448.3740344923788"
$ws.Range("F13").Value = [double]"0.350410661991486"
$ws.Range("G13").Value = [double]"6060.301438050462"
$ws.Range("H13").Value = [double]"0.5193663033388062"
$ws.Range("B14").Value = [double]"32.37400743996548"
$ws.Range("C14").Value = [double]"20"
$ws.Range("D14").Value = [double]"69"
$ws.Range("E14").Value = [double]"463.8926687782198"
$ws.Range("F14").Value = [double]"0.3283106545837414"
$ws.Range("G14").Value = [double]"6413.489705088802"
$ws.Range("H14").Value = [double]"0.5047799081095933"
$ws.Range("B15").Value = [double]"32.66068968992083"
$ws.Range("C15").Value = [double]"9.999999988282495"
$ws.Range("D15").Value = [double]"61"
$ws.Range("E15").Value = [double]"442.8869658864138"
$ws.Range("F15").Value = [double]"0.3139295491017807"
$ws.Range("G15").Value = [double]"6610.430783242102"
$ws.Range("H15").Value = [double]"0.4940780829702949"
$ws.Range("B16").Value = [double]"32.98604484994608"
$ws.Range("C16").Value = [double]"12"
$ws.Range("D16").Value = [double]"51"
$ws.Range("E16").Value = [double]"341.8566413036993"
$ws.Range("F16").Value = [double]"0.2886924322034739"
$ws.Range("G16").Value = [double]"6762.420392244242"
$ws.Range("H16").Value = [double]"0.4877845939270128"
$ws.Range("B17").Value = [double]"33.30311283997273"
$ws.Range("C17").Value = [double]"9"
$ws.Range("D17").Value = [double]"41"
$ws.Range("E17").Value = [double]"246.9476198810972"
$ws.Range("F17").Value = [double]"0.2782306195582225"
$ws.Range("G17").Value = [double]"6842.12660282243"
$ws.Range("H17").Value = [double]"0.4867362849765881"
$ws.Range("B18").Value = [double]"33.60121257994632"
$ws.Range("C18").Value = [double]"9"
$ws.Range("D18").Value = [double]"40"
$ws.Range("E18").Value = [double]"272.8626297913535"
$ws.Range("F18").Value = [double]"0.279270765184405"
$ws.Range("G18").Value = [double]"7046.135737489892"
$ws.Range("H18").Value = [double]"0.4768743298708631"
$ws.Range("B19").Value = [double]"33.85745438994953"
$ws.Range("C19").Value = [double]"8"
$ws.Range("D19").Value = [double]"43"
$ws.Range("E19").Value = [double]"287.9371385817029"
$ws.Range("F19").Value = [double]"0.2790493973473648"
$ws.Range("G19").Value = [double]"7227.972955734796"
$ws.Range("H19").Value = [double]"0.4684225383423226"
$ws.Range("B20").Value = [double]"33.90854971994773"
$ws.Range("C20").Value = [double]"7"
$ws.Range("D20").Value = [double]"26"
$ws.Range("E20").Value = [double]"163.1886755144047"
$ws.Range("F20").Value = [double]"0.2752464143539625"
$ws.Range("G20").Value = [double]"7287.317251809995"
$ws.Range("H20").Value = [double]"0.4653090917858099"
$ws.Range("B21").Value = [double]"33.95313007996212"
$ws.Range("D21").Value = [double]"32"
$ws.Range("E21").Value = [double]"247.8629601600733"
$ws.Range("F21").Value = [double]"0.2543231556529904"
$ws.Range("G21").Value = [double]"7489.864164037691"
$ws.Range("H21").Value = [double]"0.4533210394253455"
$ws.Range("B22").Value = [double]"33.96513225996133"
$ws.Range("D22").Value = [double]"19"
$ws.Range("E22").Value = [double]"144.125461568458"
$ws.Range("F22").Value = [double]"0.243230625887482"
$ws.Range("G22").Value = [double]"7541.442601454946"
$ws.Range("H22").Value = [double]"0.450379775527411"
$ws.Range("B23").Value = [double]"33.93400719996431"
$ws.Range("C23").Value = [double]"5"
$ws.Range("D23").Value = [double]"14"
$ws.Range("E23").Value = [double]"95.40279077684016"
$ws.Range("F23").Value = [double]"0.245763302454937"
$ws.Range("G23").Value = [double]"7534.768445790509"
$ws.Range("H23").Value = [double]"0.4503656276115879"
$ws.Range("B24").Value = [double]"33.86884201996345"
$ws.Range("C24").Value = [double]"14"
$ws.Range("D24").Value = [double]"31"
$ws.Range("E24").Value = [double]"185.4046056068198"
$ws.Range("F24").Value = [double]"0.2508647508589566"
$ws.Range("G24").Value = [double]"7601.683714561826"
$ws.Range("H24").Value = [double]"0.4455439517311686"
$ws.Range("B25").Value = [double]"33.66780206006461"
$ws.Range("C25").Value = [double]"54"
$ws.Range("D25").Value = [double]"92"
$ws.Range("E25").Value = [double]"410.4638908549811"
$ws.Range("F25").Value = [double]"0.2845918476479825"
$ws.Range("G25").Value = [double]"7529.976494668205"
$ws.Range("H25").Value = [double]"0.4471169609082308"
$ws.Range("B26").Value = [double]"33.45336977997385"
$ws.Range("C26").Value = [double]"48"
$ws.Range("D26").Value = [double]"148"
$ws.Range("E26").Value = [double]"603.972563241376"
$ws.Range("F26").Value = [double]"0.3089252116588242"
$ws.Range("G26").Value = [double]"7640.559580748975"
$ws.Range("H26").Value = [double]"0.4378392632950916"
$ws.Range("B27").Value = [double]"33.25321099997402"
$ws.Range("C27").Value = [double]"54.99999933210577"
$ws.Range("D27").Value = [double]"190"
$ws.Range("E27").Value = [double]"518.5679144113528"
$ws.Range("F27").Value = [double]"0.3067255742641442"
$ws.Range("G27").Value = [double]"7660.467368788697"
$ws.Range("H27").Value = [double]"0.4340885405433447"
$ws.Range("B28").Value = [double]"33.06735469998501"
$ws.Range("C28").Value = [double]"42.99999952826219"
$ws.Range("D28").Value = [double]"177"
$ws.Range("E28").Value = [double]"344.7096128525146"
$ws.Range("F28").Value = [double]"0.2970878319725"
$ws.Range("G28").Value = [double]"7796.937657950089"
$ws.Range("H28").Value = [double]"0.4241069526350276"
$ws.Range("B29").Value = [double]"32.89199017998607"
$ws.Range("C29").Value = [double]"34.0000000140303"
$ws.Range("D29").Value = [double]"155"
$ws.Range("E29").Value = [double]"196.6206573592569"
$ws.Range("F29").Value = [double]"0.2975946042943587"
$ws.Range("G29").Value = [double]"7984.585735042743"
$ws.Range("H29").Value = [double]"0.4119436032307817"
